$d = $word.ActiveDocument

# The date paragraph currently reads "03/05/2014", split across two runs:
#   run 1: "03"
#   (bookmarkStart/bookmarkEnd "_GoBack")
#   run 2: "/05/2014"
# After the edit it must read "03/06/2014", split as:
#   run 1: "03"
#   run 2: "/06"
#   (bookmarkStart/bookmarkEnd "_GoBack")
#   run 3: "/2014"

# Locate the full date string in the document and remember where it starts.
$dateRange = $d.Content
[void]$dateRange.Find.Execute("03/05/2014", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dateStart = $dateRange.Start

# Change the month from "05" to "06" (length-preserving substitution).
$monthRange = $d.Range($dateStart + 3, $dateStart + 5)
$monthRange.Text = "06"

# Relocate the "_GoBack" bookmark so it now sits right after "03/06"
# (i.e. immediately before "/2014"). Adding a bookmark with a name that
# already exists elsewhere moves it to the new location.
$bookmarkPos = $d.Range($dateStart + 5, $dateStart + 5)
[void]$d.Bookmarks.Add("_GoBack", $bookmarkPos)
